$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Datos actualizados" timestamp (A1)
$ws.Range("A1").Value = 'Datos actualizados a 29 de Marzo de 2020 a las 21:55'

# Row-level updates: country rankings shifted (new countries Colombia, Costa Rica,
# Uganda entered the table) and case statistics were refreshed for this data pull.
$rowUpdates = @{
    4 = @('Estados Unidos', 137943, 14365, 4430, 131082, 2948, 211, 2431)
    7 = @('España', 78898, 5663, 14709, 57571, 4165, 636, 6618)
    17 = @('Austria', 8743, 472, 479, 8178, 187, 18, 86)
    20 = @('Noruega', 4265, 250, 7, 4233, 91, 2, 25)
    56 = @('Colombia', 702, 94, 10, 682, 29, 4, 10)
    57 = @('Estonia', 679, 34, 20, 656, 10, 2, 3)
    58 = @('Hong Kong', 641, 81, 118, 519, 5, 0, 4)
    59 = @('Catar', 634, 44, 48, 585, 6, 0, 1)
    60 = @('Egipto', 609, 33, 132, 437, 0, 4, 40)
    77 = @('Costa Rica', 314, 19, 3, 309, 6, 0, 2)
    78 = @('Tunez', 312, 34, 2, 302, 10, 0, 8)
    79 = @('Uruguay', 304, 0, 0, 303, 9, 0, 1)
    80 = @('Taiwan', 298, 15, 39, 257, 0, 0, 2)
    81 = @('Kazajistan', 284, 56, 20, 263, 0, 0, 1)
    82 = @('Moldavia', 263, 32, 2, 259, 33, 0, 2)
    83 = @('Republica de Macedonia', 259, 18, 3, 250, 1, 2, 6)
    84 = @('Jordania', 259, 13, 18, 239, 3, 1, 2)
    85 = @('Kuwait', 255, 20, 67, 188, 12, 0, 0)
    134 = @('Uganda', 33, 3, 0, 33, 0, 0, 0)
    135 = @('Guam', 32, 0, 0, 31, 0, 0, 1)
    136 = @('Jamaica', 32, 2, 2, 32, 0, 0, 1)
    137 = @('Polinesia Francesa', 30, 0, 0, 30, 0, 0, 0)
}

foreach ($rowNum in $rowUpdates.Keys) {
    $values = $rowUpdates[$rowNum]
    $ws.Cells.Item([int]$rowNum, 1).Value = $values[0]
    for ($col = 1; $col -lt $values.Length; $col++) {
        $ws.Cells.Item([int]$rowNum, $col + 1).Value = $values[$col]
    }
}
